$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the STORE column data (column L) for all data rows - values were "PASS"
$ws.Range("L2:L11").ClearContents()

# Clear empty placeholder cells that are no longer needed
$ws.Range("F3:F11").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("I6:I11").ClearContents()
$ws.Range("K3:K11").ClearContents()
